# Auto-generated edit script for RNN Results.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  104 = @{ F='python rnn.py NSAA position all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  105 = @{ F='python rnn.py NSAA position all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  106 = @{ F='python rnn.py NSAA position all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  107 = @{ F='python rnn.py NSAA sensorMagneticField all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  108 = @{ F='python rnn.py NSAA sensorMagneticField all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  109 = @{ F='python rnn.py NSAA sensorMagneticField all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  110 = @{ F='python rnn.py NSAA jointAngle all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  111 = @{ F='python rnn.py NSAA jointAngle all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  112 = @{ F='python rnn.py NSAA jointAngle all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  113 = @{ F='python rnn.py NSAA jointAngleXZY all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  114 = @{ F='python rnn.py NSAA jointAngleXZY all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  115 = @{ C='Raw joint angles from allmatfiles to perform overall NSAA score regression'; D='python ext_raw_measures.py allmatfiles all jointAngle'; F='python rnn.py NSAA jointAngleXZY all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9'; G='x'; style=$false }
  116 = @{ F='python rnn.py NSAA AD all dhc --seq_len=10 --seq_overlap=0.9 --epochs=300'; G='x'; style=$false }
  117 = @{ F='python rnn.py NSAA AD all overall --seq_len=10 --seq_overlap=0.9 --epochs=300'; G='x'; style=$false }
  118 = @{ F='python rnn.py NSAA AD all acts --seq_len=10 --seq_overlap=0.9 --epochs=300'; G='x'; style=$false }
  121 = @{ C='NOTE: received many more files (15/06), hence more available training data'; style=$false }
  122 = @{ F='python ext_raw_measures.py 6minwalk-matfiles all all'; G='x'; style=$true }
  123 = @{ C='(NOT SURE IF POSSIBLE)'; D='Single-act stat values from NSAA\AD to perform D/HC classification'; E='python mat_act_div.py V1 all; python comp_stat_vals.py NSAA AD all --split_size=1  --single_act'; F='python ext_raw_measures.py 6MW-matFiles all all'; G='x'; style=$true }
  124 = @{ F='python comp_stat_vals.py 6MW-matFiles AD all --split_size=1'; G='x'; style=$true }
  125 = @{ F='python rnn.py 6minwalk-matfiles position all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  126 = @{ F='python rnn.py 6minwalk-matfiles position all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  127 = @{ F='python rnn.py 6minwalk-matfiles position all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  128 = @{ F='python rnn.py 6minwalk-matfiles sensorMagneticField all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  129 = @{ F='python rnn.py 6minwalk-matfiles sensorMagneticField all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  130 = @{ F='python rnn.py 6minwalk-matfiles sensorMagneticField all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  131 = @{ F='python rnn.py 6minwalk-matfiles jointAngle all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  132 = @{ F='python rnn.py 6minwalk-matfiles jointAngle all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  133 = @{ F='python rnn.py 6minwalk-matfiles jointAngle all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  134 = @{ F='python rnn.py 6minwalk-matfiles jointAngleXZY all dhc --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  135 = @{ F='python rnn.py 6minwalk-matfiles jointAngleXZY all overall --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  136 = @{ F='python rnn.py 6minwalk-matfiles jointAngleXZY all acts --seq_len=600 --seq_overlap=0.9 --discard_prop=0.9 --other_dir=6MW-matFiles'; G='x'; style=$false }
  137 = @{ F='python ft_sel_red.py 6MW-matFiles AD all pca --num_features=30 --no_normalize'; G='x'; style=$false }
  138 = @{ F='python rnn.py 6minwalk-matfiles AD all dhc --seq_len=10 --seq_overlap=0.9 --epochs=300 --other_dir=6MW-matFiles'; G='x'; style=$false }
  139 = @{ F='python rnn.py 6minwalk-matfiles AD all overall --seq_len=10 --seq_overlap=0.9 --epochs=300 --other_dir=6MW-matFiles'; G='x'; style=$false }
  140 = @{ F='python rnn.py 6minwalk-matfiles AD all acts --seq_len=10 --seq_overlap=0.9 --epochs=300 --other_dir=6MW-matFiles'; G='x'; style=$false }
}

foreach ($rowNum in ($rows.Keys | Sort-Object)) {
  $cellData = $rows[$rowNum]
  foreach ($col in @('C','D','E','F','G')) {
    if ($cellData.ContainsKey($col)) {
      $ws.Range("$col$rowNum").Value = $cellData[$col]
    }
  }
}

# Apply the distinct font/alignment style (fontId=1, black color, vertical-center)
# to F122 directly, then propagate via copy/paste-format so only a single new
# cellXf entry is created (matches the target styles.xml: cellXfs count=2).
$ws.Range("F122").Font.Color = 0
$ws.Range("F122").VerticalAlignment = -4108
$ws.Range("F122").Copy()
$ws.Range("F123:F124").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection to match the saved workbook state
$ws.Range("H140").Select()

